$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the new "Untouched" worksheet (after "Tables") with its sample data.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Untouched"
$ws2.Range("A1").Value = "some values"
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5

# ---------------------------------------------------------------------------
# 2. New defined name "Formula" that references the helper sheet.
# ---------------------------------------------------------------------------
$wb.Names.Add("Formula", "=IF(Tables!`$G`$7>1,Tables!`$G`$6,SUM(Untouched!`$A`$2:`$A`$6)/COUNT(Untouched!`$A`$2:`$A`$6))")

# ---------------------------------------------------------------------------
# 3. New helper cells on the "Tables" sheet (columns G/H).
# ---------------------------------------------------------------------------
$ws1.Range("B2:C2").Copy()
$ws1.Range("G2:H2").PasteSpecial(-4122)
$ws1.Range("G2").Value = "defined name with formula"
$ws1.Range("G2:H2").Merge()

$ws1.Range("G3").Value = "helper"
$ws1.Range("H3").Value = "formula"

$ws1.Range("G6").Value = 0
$ws1.Range("G7").Value = 1

$ws1.Range("H5").FormulaArray = "=Formula"

$ws1.Columns.Item(7).ColumnWidth = 16.95
$ws1.Columns.Item(8).ColumnWidth = 12.75

# ---------------------------------------------------------------------------
# 4. Clear the stray border formatting on C14:D14.
# ---------------------------------------------------------------------------
$ws1.Range("C14:D14").ClearFormats()

# ---------------------------------------------------------------------------
# 5. Re-localise the built-in cell styles (German Excel names them
#    differently) - best-effort, mirrors the author's environment.
# ---------------------------------------------------------------------------
try { $wb.Styles.Item("Heading 1").Name = "Überschrift 1" } catch {}
try { $wb.Styles.Item("Normal").Name = "Standard" } catch {}
try { $wb.Styles.Item("Percent").Name = "Prozent" } catch {}

# ---------------------------------------------------------------------------
# 6. Restore the selection/view state on the "Tables" sheet.
# ---------------------------------------------------------------------------
$ws1.Select() | Out-Null
$ws1.Range("G8").Select() | Out-Null
